# Applies the "fix issues de compilation" edit:
#   - rename the "Include from Medication types" sheet to "Include #0"
#   - bump the Metadata "Date" value
#   - insert a new "Jurisdiction" property row right after "Contact",
#     pushing Description / Purpose / Copyright / Immutable down by one row

$wb = $excel.ActiveWorkbook

# --- 1) Rename the include sheet -------------------------------------------------
$wsInclude = $wb.Worksheets.Item("Include from Medication types")
$wsInclude.Name = "Include #0"

# --- 2) Edit the Metadata sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("Metadata")

# Give the about-to-be-created row 15 the same look (borders/font/wrap) as row 14
# before we start moving values around.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

# Shift rows 11-14 (Description, Purpose, Copyright, Immutable) down to 12-15,
# working from the bottom up so we never overwrite a value before reading it.
$ws.Cells.Item(15, 1).Value = $ws.Cells.Item(14, 1).Value2
$ws.Cells.Item(15, 2).Value = $ws.Cells.Item(14, 2).Value2

$ws.Cells.Item(14, 1).Value = $ws.Cells.Item(13, 1).Value2
$ws.Cells.Item(14, 2).Value = $ws.Cells.Item(13, 2).Value2

$ws.Cells.Item(13, 1).Value = $ws.Cells.Item(12, 1).Value2
$ws.Cells.Item(13, 2).Value = $ws.Cells.Item(12, 2).Value2

$ws.Cells.Item(12, 1).Value = $ws.Cells.Item(11, 1).Value2
$ws.Cells.Item(12, 2).Value = $ws.Cells.Item(11, 2).Value2

# Row 11 becomes the new "Jurisdiction" property (empty value).
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Bump the publication Date.
$ws.Cells.Item(8, 2).Value = "2024-09-12T14:01:50+00:00"
